$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 299, shifting existing rows 299-380 down to 300-381.
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new record's data.
$ws.Cells.Item(299, 1).Value = 8
$ws.Cells.Item(299, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(299, 3).Value = "Coquimbo"
$ws.Cells.Item(299, 4).Value = 44642
$ws.Cells.Item(299, 5).Value = 4
$ws.Cells.Item(299, 6).Value = 100114001
$ws.Cells.Item(299, 7).Value = "Papa"
$ws.Cells.Item(299, 8).Value = "Asterix"
$ws.Cells.Item(299, 9).Value = "1a (cosecha)"
$ws.Cells.Item(299, 10).Value = 2000
$ws.Cells.Item(299, 11).Value = 8500
$ws.Cells.Item(299, 12).Value = 9000
$ws.Cells.Item(299, 13).Value = 8750
$ws.Cells.Item(299, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(299, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(299, 16).Value = 350
$ws.Cells.Item(299, 17).Value = 25
$ws.Cells.Item(299, 18).Value = "Hortaliza"
